$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1950
$ws.Cells.Item(6, 9).Value = 100
$ws.Cells.Item(6, 10).Value = 3800
$ws.Cells.Item(6, 11).Value = 300
$ws.Cells.Item(6, 12).Value = 11400
$ws.Cells.Item(6, 13).Value = -188
$ws.Cells.Item(6, 14).Value = -11624

$ws.Cells.Item(7, 8).Value = 34998.5
$ws.Cells.Item(7, 10).Value = 34998.5
$ws.Cells.Item(7, 12).Value = 34998.5
$ws.Cells.Item(7, 14).Value = -35222.5

$ws.Cells.Item(10, 8).Value = 29999
$ws.Cells.Item(10, 10).Value = 29999
$ws.Cells.Item(10, 12).Value = 29999
$ws.Cells.Item(10, 14).Value = -30585

$ws.Cells.Item(14, 8).Value = 34998.5
$ws.Cells.Item(14, 10).Value = 34998.5
$ws.Cells.Item(14, 12).Value = 34998.5
$ws.Cells.Item(14, 14).Value = -35380.5

$ws.Cells.Item(16, 8).Value = 35000
$ws.Cells.Item(16, 10).Value = 35000
$ws.Cells.Item(16, 12).Value = 35000
$ws.Cells.Item(16, 14).Value = -35460

$ws.Cells.Item(19, 8).Value = 1144940.9
$ws.Cells.Item(19, 9).Value = 2393038.2
$ws.Cells.Item(19, 10).Value = 851.5
$ws.Cells.Item(19, 11).Value = 2393038.2
$ws.Cells.Item(19, 12).Value = 851.5
$ws.Cells.Item(19, 13).Value = -2392863.2
$ws.Cells.Item(19, 14).Value = -1201.5

$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 14).ClearContents()

$ws.Cells.Item(55, 8).Value = 176.63637
$ws.Cells.Item(55, 9).Value = 133.85715
$ws.Cells.Item(55, 10).Value = 196.6
$ws.Cells.Item(55, 11).Value = 133.85715
$ws.Cells.Item(55, 12).Value = 196.6
$ws.Cells.Item(55, 13).Value = 80.14285000000001
$ws.Cells.Item(55, 14).Value = -624.6

$ws.Cells.Item(111, 8).Value = 1079.8572
$ws.Cells.Item(111, 9).Value = 1500
$ws.Cells.Item(111, 11).Value = 4500
$ws.Cells.Item(111, 13).Value = -1433

$ws.Cells.Item(115, 8).Value = 1502.7778
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 14).ClearContents()

$ws.Cells.Item(131, 8).Value = 2373.7222
$ws.Cells.Item(131, 9).Value = 504.625
$ws.Cells.Item(131, 10).Value = 3869
$ws.Cells.Item(131, 11).Value = 1513.875
$ws.Cells.Item(131, 12).Value = 11607
$ws.Cells.Item(131, 13).Value = 3526.125
$ws.Cells.Item(131, 14).Value = -21687

$ws.Cells.Item(132, 8).Value = 327513.06
$ws.Cells.Item(132, 9).Value = 4845.643
$ws.Cells.Item(132, 11).Value = 14536.929
$ws.Cells.Item(132, 13).Value = -12006.929

$ws.Cells.Item(138, 8).Value = 5569.59
$ws.Cells.Item(138, 9).Value = 776.875
$ws.Cells.Item(138, 10).Value = 7083.079
$ws.Cells.Item(138, 11).Value = 2330.625
$ws.Cells.Item(138, 12).Value = 21249.237
$ws.Cells.Item(138, 13).Value = 2809.375
$ws.Cells.Item(138, 14).Value = -31529.237

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 4316.2666
$ws.Cells.Item(74, 9).Value = 4599.524
$ws.Cells.Item(74, 11).Value = 4599.524
$ws.Cells.Item(74, 13).Value = -3725.524

$ws.Cells.Item(77, 8).Value = 4316.2666
$ws.Cells.Item(77, 9).Value = 4599.524
$ws.Cells.Item(77, 11).Value = 22997.62
$ws.Cells.Item(77, 13).Value = -18629.62

$ws.Cells.Item(122, 8).Value = 1859.4
$ws.Cells.Item(122, 9).Value = 855.2857
$ws.Cells.Item(122, 11).Value = 2565.8571
$ws.Cells.Item(122, 13).Value = -115.8571000000002

$ws.Cells.Item(137, 8).Value = 44546
$ws.Cells.Item(137, 10).Value = 44546
$ws.Cells.Item(137, 12).Value = 44546
$ws.Cells.Item(137, 14).Value = -54746

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 3349258
$ws.Cells.Item(7, 9).Value = 4010379.2
$ws.Cells.Item(7, 10).Value = 3094980.8
$ws.Cells.Item(7, 11).Value = 4010379.2
$ws.Cells.Item(7, 12).Value = 3094980.8
$ws.Cells.Item(7, 13).Value = -4010266.2
$ws.Cells.Item(7, 14).Value = -3095206.8

$ws.Cells.Item(99, 8).Value = 2826.4092
$ws.Cells.Item(99, 9).Value = 1055.7858
$ws.Cells.Item(99, 11).Value = 1055.7858
$ws.Cells.Item(99, 13).Value = 442.2141999999999

$ws.Cells.Item(105, 8).Value = 5292632.5
$ws.Cells.Item(105, 9).Value = 5557184
$ws.Cells.Item(105, 10).Value = 1598
$ws.Cells.Item(105, 11).Value = 5557184
$ws.Cells.Item(105, 12).Value = 1598
$ws.Cells.Item(105, 13).Value = -5555437
$ws.Cells.Item(105, 14).Value = -5092

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 1900.8096
$ws.Cells.Item(105, 9).Value = 1971.3846
$ws.Cells.Item(105, 10).Value = 1786.125
$ws.Cells.Item(105, 11).Value = 1971.3846
$ws.Cells.Item(105, 12).Value = 1786.125
$ws.Cells.Item(105, 13).Value = -224.3846000000001
$ws.Cells.Item(105, 14).Value = -5280.125

$ws.Cells.Item(107, 8).Value = 753.45
$ws.Cells.Item(107, 9).Value = 689.4167
$ws.Cells.Item(107, 10).Value = 849.5
$ws.Cells.Item(107, 11).Value = 689.4167
$ws.Cells.Item(107, 12).Value = 849.5
$ws.Cells.Item(107, 13).Value = 1230.5833
$ws.Cells.Item(107, 14).Value = -4689.5

$ws.Cells.Item(122, 8).Value = 2313.9
$ws.Cells.Item(122, 9).Value = 1773
$ws.Cells.Item(122, 11).Value = 5319
$ws.Cells.Item(122, 13).Value = -2869

$ws.Cells.Item(132, 8).Value = 3254.818
$ws.Cells.Item(132, 9).Value = 1755.8182
$ws.Cells.Item(132, 10).Value = 4753.8184
$ws.Cells.Item(132, 11).Value = 5267.4546
$ws.Cells.Item(132, 12).Value = 14261.4552
$ws.Cells.Item(132, 13).Value = -2737.4546
$ws.Cells.Item(132, 14).Value = -19321.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 1483.28
$ws.Cells.Item(3, 9).Value = 1031
$ws.Cells.Item(3, 11).Value = 3093
$ws.Cells.Item(3, 13).Value = -2981

$ws.Cells.Item(7, 8).Value = 149.5
$ws.Cells.Item(7, 9).Value = 99.40000000000001
$ws.Cells.Item(7, 10).Value = 400
$ws.Cells.Item(7, 11).Value = 298.2
$ws.Cells.Item(7, 12).Value = 1200
$ws.Cells.Item(7, 13).Value = -186.2
$ws.Cells.Item(7, 14).Value = -1424

$ws.Cells.Item(10, 8).Value = 416.25
$ws.Cells.Item(10, 9).Value = 332.85715
$ws.Cells.Item(10, 10).Value = 1000
$ws.Cells.Item(10, 11).Value = 998.5714499999999
$ws.Cells.Item(10, 12).Value = 3000
$ws.Cells.Item(10, 13).Value = -859.5714499999999
$ws.Cells.Item(10, 14).Value = -3278

$ws.Cells.Item(104, 8).Value = 1999.6666
$ws.Cells.Item(104, 10).Value = 1999.6666
$ws.Cells.Item(104, 12).Value = 5998.9998
$ws.Cells.Item(104, 14).Value = -11240.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 2504356.2
$ws.Cells.Item(3, 9).Value = 5002068.5
$ws.Cells.Item(3, 10).Value = 6643.857
$ws.Cells.Item(3, 11).Value = 5002068.5
$ws.Cells.Item(3, 12).Value = 6643.857
$ws.Cells.Item(3, 13).Value = -5001952.5
$ws.Cells.Item(3, 14).Value = -6875.857

$ws.Cells.Item(10, 8).Value = 8184618
$ws.Cells.Item(10, 9).Value = 10001200
$ws.Cells.Item(10, 10).Value = 9998.5
$ws.Cells.Item(10, 11).Value = 10001200
$ws.Cells.Item(10, 12).Value = 9998.5
$ws.Cells.Item(10, 13).Value = -10001031
$ws.Cells.Item(10, 14).Value = -10336.5

$ws.Cells.Item(11, 8).Value = 10940647
$ws.Cells.Item(11, 10).Value = 1805034.8
$ws.Cells.Item(11, 12).Value = 1805034.8
$ws.Cells.Item(11, 14).Value = -1805312.8

$ws.Cells.Item(46, 8).Value = 31744.857
$ws.Cells.Item(46, 10).Value = 31744.857
$ws.Cells.Item(46, 12).Value = 31744.857
$ws.Cells.Item(46, 14).Value = -32056.857

$ws.Cells.Item(126, 8).Value = 3422.9167
$ws.Cells.Item(126, 9).Value = 2975.6758
$ws.Cells.Item(126, 11).Value = 8927.027399999999
$ws.Cells.Item(126, 13).Value = -6457.027399999999

$ws.Cells.Item(132, 8).Value = 4415.7144
$ws.Cells.Item(132, 9).Value = 978
$ws.Cells.Item(132, 10).Value = 8999.333000000001
$ws.Cells.Item(132, 11).Value = 2934
$ws.Cells.Item(132, 12).Value = 26997.999
$ws.Cells.Item(132, 13).Value = -404
$ws.Cells.Item(132, 14).Value = -32057.999

$ws.Cells.Item(137, 8).Value = 39830
$ws.Cells.Item(137, 10).Value = 39830
$ws.Cells.Item(137, 12).Value = 39830
$ws.Cells.Item(137, 14).Value = -50030

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2211.9
$ws.Cells.Item(61, 9).Value = 2400
$ws.Cells.Item(61, 10).Value = 2086.5
$ws.Cells.Item(61, 11).Value = 2400
$ws.Cells.Item(61, 12).Value = 2086.5
$ws.Cells.Item(61, 13).Value = -2198
$ws.Cells.Item(61, 14).Value = -2490.5

$ws.Cells.Item(107, 8).Value = 1495
$ws.Cells.Item(107, 9).Value = 1495
$ws.Cells.Item(107, 11).Value = 1495
$ws.Cells.Item(107, 13).Value = 425

$ws.Cells.Item(113, 8).Value = 2211.9
$ws.Cells.Item(113, 9).Value = 2400
$ws.Cells.Item(113, 10).Value = 2086.5
$ws.Cells.Item(113, 11).Value = 2400
$ws.Cells.Item(113, 12).Value = 2086.5
$ws.Cells.Item(113, 13).Value = -230
$ws.Cells.Item(113, 14).Value = -6426.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 13335913
$ws.Cells.Item(132, 9).Value = 1559.0667
$ws.Cells.Item(132, 10).Value = 33337444
$ws.Cells.Item(132, 11).Value = 4677.2001
$ws.Cells.Item(132, 12).Value = 100012332
$ws.Cells.Item(132, 13).Value = -2147.2001
$ws.Cells.Item(132, 14).Value = -100017392

$ws.Cells.Item(136, 8).Value = 4058.8
$ws.Cells.Item(136, 9).Value = 2345.4285
$ws.Cells.Item(136, 11).Value = 7036.2855
$ws.Cells.Item(136, 13).Value = -4486.2855
